$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.084.56'
$ws.Cells.Item(2, 5).Value = '  -0.24%  '
$ws.Cells.Item(3, 4).Value = '1.757.13'
$ws.Cells.Item(3, 5).Value = '  -2.74%  '
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(5, 4).Value = '337.76'
$ws.Cells.Item(5, 5).Value = '  -0.52%  '
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
$ws.Cells.Item(7, 4).Value = '0.3776'
$ws.Cells.Item(7, 5).Value = '  -4.18%  '
$ws.Cells.Item(8, 4).Value = '0.3350'
$ws.Cells.Item(8, 5).Value = '  -4.09%  '
$ws.Cells.Item(9, 4).Value = '45.33'
$ws.Cells.Item(9, 5).Value = '  -6.05%  '
$ws.Cells.Item(10, 4).Value = '1.119'
$ws.Cells.Item(10, 5).Value = '  -4.42%  '
$ws.Cells.Item(11, 4).Value = '0.07207'
$ws.Cells.Item(11, 5).Value = '  -4.40%  '
$ws.Cells.Item(12, 4).Value = '1.001'
$ws.Cells.Item(12, 5).Value = '  -0.11%  '
$ws.Cells.Item(13, 4).Value = '22.61'
$ws.Cells.Item(13, 5).Value = '  +2.63%  '
$ws.Cells.Item(14, 4).Value = '6.148'
$ws.Cells.Item(14, 5).Value = '  -5.55%  '
$ws.Cells.Item(15, 4).Value = '7.169'
$ws.Cells.Item(15, 5).Value = '  +0.31%  '
$ws.Cells.Item(16, 4).Value = '1.758.98'
$ws.Cells.Item(16, 5).Value = '  -2.70%  '
$ws.Cells.Item(17, 5).Value = '  -4.04%  '
$ws.Cells.Item(18, 4).Value = '0.06602'
$ws.Cells.Item(18, 5).Value = '  -1.66%  '
$ws.Cells.Item(19, 4).Value = '80.78'
$ws.Cells.Item(19, 5).Value = '  -5.00%  '
$ws.Cells.Item(20, 4).Value = '1.001'
$ws.Cells.Item(20, 5).Value = '  -0.07%  '
$ws.Cells.Item(21, 4).Value = '16.94'
$ws.Cells.Item(21, 5).Value = '  -4.40%  '
$ws.Cells.Item(22, 4).Value = '6.243'
$ws.Cells.Item(22, 5).Value = '  -4.75%  '
$ws.Cells.Item(23, 4).Value = '28.069.27'
$ws.Cells.Item(23, 5).Value = '  -0.30%  '
$ws.Cells.Item(24, 4).Value = '11.67'
$ws.Cells.Item(24, 5).Value = '  -5.89%  '
$ws.Cells.Item(25, 4).Value = '2.404'
$ws.Cells.Item(25, 5).Value = '  -0.02%  '
$ws.Cells.Item(26, 4).Value = '152.95'
$ws.Cells.Item(26, 5).Value = '  -0.83%  '
$ws.Cells.Item(27, 4).Value = '19.89'
$ws.Cells.Item(27, 5).Value = '  -7.08%  '
$ws.Cells.Item(28, 4).Value = '2.328'
$ws.Cells.Item(28, 5).Value = '  -7.43%  '
$ws.Cells.Item(29, 4).Value = '1.957.27'
$ws.Cells.Item(29, 5).Value = '  -2.79%  '
$ws.Cells.Item(30, 4).Value = '131.89'
$ws.Cells.Item(30, 5).Value = '  -2.58%  '
$ws.Cells.Item(31, 4).Value = '1.249'
$ws.Cells.Item(31, 5).Value = '  -16.02%  '
$ws.Cells.Item(32, 4).Value = '4.024'
$ws.Cells.Item(32, 5).Value = '  +0.17%  '
$ws.Cells.Item(33, 4).Value = '5.786'
$ws.Cells.Item(33, 5).Value = '  -6.49%  '
$ws.Cells.Item(34, 4).Value = '0.08762'
$ws.Cells.Item(34, 5).Value = '  -0.91%  '
$ws.Cells.Item(35, 4).Value = '12.24'
$ws.Cells.Item(35, 5).Value = '  -6.62%  '
$ws.Cells.Item(36, 4).Value = '0.02337'
$ws.Cells.Item(36, 5).Value = '  -3.69%  '
$ws.Cells.Item(37, 4).Value = '0.6671'
$ws.Cells.Item(37, 5).Value = '  -3.84%  '
$ws.Cells.Item(38, 4).Value = '0.06195'
$ws.Cells.Item(38, 5).Value = '  -5.13%  '
$ws.Cells.Item(39, 4).Value = '5.170'
$ws.Cells.Item(39, 5).Value = '  -5.11%  '
$ws.Cells.Item(40, 4).Value = '0.2111'
$ws.Cells.Item(40, 5).Value = '  -4.52%  '
$ws.Cells.Item(41, 4).Value = '1.216'
$ws.Cells.Item(41, 5).Value = '  -3.04%  '
$ws.Cells.Item(42, 4).Value = '1.443'
$ws.Cells.Item(42, 5).Value = '  -10.29%  '
$ws.Cells.Item(43, 4).Value = '8.024'
$ws.Cells.Item(43, 5).Value = '  -5.73%  '
$ws.Cells.Item(44, 5).Value = '  +0.01%  '
$ws.Cells.Item(45, 4).Value = '13.64'
$ws.Cells.Item(45, 5).Value = '  -6.85%  '
$ws.Cells.Item(46, 2).Value = 'PancakeSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(46, 4).Value = '3.835'
$ws.Cells.Item(46, 5).Value = '  -1.04%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '0.6042'
$ws.Cells.Item(47, 5).Value = '  -5.88%  '
$ws.Cells.Item(48, 4).Value = '129.08'
$ws.Cells.Item(48, 5).Value = '  -1.41%  '
$ws.Cells.Item(49, 4).Value = '2.022'
$ws.Cells.Item(49, 5).Value = '  -5.69%  '
$ws.Cells.Item(50, 4).Value = '1.182'
$ws.Cells.Item(50, 5).Value = '  +1.83%  '
$ws.Cells.Item(51, 4).Value = '0.07179'
$ws.Cells.Item(51, 5).Value = '  -0.17%  '
